$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.Formula = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "306.22"
Set-TextValue "E2" "-0.74%"
Set-TextValue "D3" "39.33"
Set-TextValue "E3" "8.11%"
Set-TextValue "D4" "5.115"
Set-TextValue "E4" "1.26%"
Set-TextValue "D5" "0.08075"
Set-TextValue "E5" "-0.68%"
Set-TextValue "D6" "1.928"
Set-TextValue "E6" "-3.24%"
Set-TextValue "D7" "4.200"
Set-TextValue "E7" "1.21%"
Set-TextValue "D8" "8.056"
Set-TextValue "D9" "0.9259"
Set-TextValue "E9" "-0.17%"
Set-TextValue "D10" "0.1393"
Set-TextValue "E10" "-5.69%"
Set-TextValue "D11" "0.1917"
Set-TextValue "E11" "-1.04%"
Set-TextValue "D12" "0.08990"
Set-TextValue "E12" "-1.09%"
Set-TextValue "D13" "0.03523"
Set-TextValue "E13" "-0.10%"
Set-TextValue "D14" "0.09779"
Set-TextValue "E14" "-0.90%"
Set-TextValue "D15" "0.001394"
Set-TextValue "E15" "-1.04%"
Set-TextValue "D16" "0.005857"
Set-TextValue "E16" "-10.56%"
Set-TextValue "D17" "3.764"
Set-TextValue "E17" "-2.18%"
Set-TextValue "E18" "-1.27%"
Set-TextValue "D20" "0.1293"
Set-TextValue "E20" "-1.49%"
Set-TextValue "D21" "4.686"
Set-TextValue "E21" "-2.59%"
Set-TextValue "E22" "3.06%"
Set-TextValue "D23" "0.04371"
Set-TextValue "E23" "0.08%"
Set-TextValue "D24" "0.001205"
Set-TextValue "E24" "-2.39%"
Set-TextValue "D25" "0.004280"
Set-TextValue "E25" "2.85%"
Set-TextValue "E26" "0.04%"
Set-TextValue "D39" "0.02036"
Set-TextValue "E39" "-4.45%"
Set-TextValue "D40" "0.05028"
Set-TextValue "E40" "-1.77%"
Set-TextValue "D41" "0.007530"
Set-TextValue "E41" "0.84%"
Set-TextValue "D42" "0.009701"
Set-TextValue "E42" "-3.67%"
Set-TextValue "D43" "0.1343"
Set-TextValue "E43" "-1.88%"
Set-TextValue "E44" "-1.84%"
Set-TextValue "D45" "0.009795"
Set-TextValue "E45" "0.96%"
Set-TextValue "D46" "0.00006205"
Set-TextValue "E46" "-1.07%"
Set-TextValue "E47" "0.02%"
Set-TextValue "D48" "0.002786"
Set-TextValue "E49" "12.58%"
Set-TextValue "E50" "0.02%"
Set-TextValue "E51" "0.02%"
